# Applies the "Updated symbol list" data refresh described in the commit:
# - Coin/Link columns (B:C) for rows 7-19 rotate up by one row (the GateToken row
#   wraps from row 7 down to row 19), reflecting the source list re-ranking.
# - Price (D) and Volume(1h) (E) columns get refreshed quote text for every affected row.
#
# All of these columns are stored as literal text in the sheet (e.g. "289.10", "1.22%"),
# so each numeric-looking cell has its NumberFormat forced to Text ("@") right before the
# value is written -- otherwise Excel would silently reinterpret the text as a number/
# percentage (losing the exact formatting, e.g. "289.10" -> 289.1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin (B) / Link (C) columns: rows 7-19 rotate up by one row ------------------------
$ws.Range('B7').Value = 'FTXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('B15').Value = 'CoinExToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('B16').Value = 'One'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('B19').Value = 'GateToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'

# --- Price (D) / Volume 1h (E) columns: refreshed text values -----------------------------
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '289.10'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '1.22%'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '29.35'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '1.90%'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '5.100'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '3.29%'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.06677'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '2.76%'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '7.368'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '1.87%'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.358'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '0.86%'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.9166'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '0.46%'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.1588'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '2.63%'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.06678'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '1.45%'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07711'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '-0.24%'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.02934'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '-1.38%'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.08992'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '0.09%'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.001572'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '-1.77%'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.04514'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '0.64%'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0006457'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '-1.30%'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.006276'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '3.93%'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.449'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '-0.30%'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.408'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '0.49%'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '-0.89%'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '2.05%'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '-2.97%'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.061'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '0.47%'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.1568'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '0.55%'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.001189'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '-0.28%'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.004127'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '-4.44%'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '5.54%'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '-1.41%'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '1.49%'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.006748'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '1.05%'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.1239'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '0.49%'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.01333'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '13.26%'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '-9.52%'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.00005711'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '4.15%'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.972'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '26.20%'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '-29.59%'
